# Rewrite the "12. Training" section.
#
# The section currently consists of four placeholder paragraphs that each
# just say "Sample ", followed by an (otherwise empty) paragraph that carries
# the _GoBack bookmark. We need to:
#   1. Remove the four "Sample " placeholder paragraphs.
#   2. Turn the remaining (bookmark) paragraph into the real training blurb,
#      first-line indented by half an inch, while keeping the bookmark in
#      place.

$d = $word.ActiveDocument

# Find every paragraph whose entire text is the placeholder "Sample " (plus
# its trailing paragraph mark) so we don't have to hard-code paragraph
# indices.
$sampleIndexes = New-Object System.Collections.ArrayList
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Sample `r") {
        [void]$sampleIndexes.Add($i)
    }
    $i = $i + 1
}

# Delete them back-to-front so earlier indices stay valid while we work.
for ($j = $sampleIndexes.Count - 1; $j -ge 0; $j--) {
    $d.Paragraphs.Item($sampleIndexes[$j]).Range.Delete()
}

# The paragraph that used to follow the placeholders (and holds the
# _GoBack bookmark) now occupies the slot where the first placeholder was.
$targetIndex = $sampleIndexes[0]

# Insert the real training copy right after the (zero-length) bookmark so
# the bookmark still wraps no text, matching the original layout.
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertAfter("15 hours personal training will be provided from the environment familiarization up to actual usage and control.")

# Give the paragraph a half-inch first-line indent.
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.ParagraphFormat.FirstLineIndent = $word.InchesToPoints(0.5)
